# Auto-generated edit script
# Updates column F (想去人数 / 'want to go' count) across 4 worksheets
# matching the commit 'Update gh-pages to output generated at 456a3b4'

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 1826
$ws.Cells.Item(4, 6).Value = 22
$ws.Cells.Item(5, 6).Value = 34
$ws.Cells.Item(6, 6).Value = 1076
$ws.Cells.Item(7, 6).Value = 40
$ws.Cells.Item(8, 6).Value = 167
$ws.Cells.Item(9, 6).Value = 566
$ws.Cells.Item(10, 6).Value = 50
$ws.Cells.Item(11, 6).Value = 444
$ws.Cells.Item(12, 6).Value = 203
$ws.Cells.Item(13, 6).Value = 1358
$ws.Cells.Item(14, 6).Value = 1193
$ws.Cells.Item(15, 6).Value = 1395
$ws.Cells.Item(16, 6).Value = 12
$ws.Cells.Item(17, 6).Value = 109
$ws.Cells.Item(18, 6).Value = 274
$ws.Cells.Item(19, 6).Value = 1539
$ws.Cells.Item(21, 6).Value = 763
$ws.Cells.Item(22, 6).Value = 307
$ws.Cells.Item(23, 6).Value = 44
$ws.Cells.Item(24, 6).Value = 105
$ws.Cells.Item(25, 6).Value = 1167
$ws.Cells.Item(26, 6).Value = 307
$ws.Cells.Item(27, 6).Value = 21
$ws.Cells.Item(29, 6).Value = 541
$ws.Cells.Item(30, 6).Value = 995
$ws.Cells.Item(31, 6).Value = 213266
$ws.Cells.Item(35, 6).Value = 875
$ws.Cells.Item(36, 6).Value = 16
$ws.Cells.Item(37, 6).Value = 15
$ws.Cells.Item(38, 6).Value = 815
$ws.Cells.Item(39, 6).Value = 1539
$ws.Cells.Item(40, 6).Value = 83
$ws.Cells.Item(41, 6).Value = 13
$ws.Cells.Item(42, 6).Value = 774
$ws.Cells.Item(44, 6).Value = 758
$ws.Cells.Item(45, 6).Value = 100

# --- Sheet 2: 演出 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value = 106
$ws.Cells.Item(8, 6).Value = 83
$ws.Cells.Item(11, 6).Value = 1371
$ws.Cells.Item(13, 6).Value = 2460
$ws.Cells.Item(14, 6).Value = 1164
$ws.Cells.Item(15, 6).Value = 381
$ws.Cells.Item(16, 6).Value = 710
$ws.Cells.Item(17, 6).Value = 199
$ws.Cells.Item(19, 6).Value = 59
$ws.Cells.Item(20, 6).Value = 14
$ws.Cells.Item(22, 6).Value = 414
$ws.Cells.Item(23, 6).Value = 17
$ws.Cells.Item(24, 6).Value = 6
$ws.Cells.Item(25, 6).Value = 267
$ws.Cells.Item(26, 6).Value = 44140
$ws.Cells.Item(27, 6).Value = 10
$ws.Cells.Item(31, 6).Value = 227
$ws.Cells.Item(33, 6).Value = 45
$ws.Cells.Item(38, 6).Value = 161
$ws.Cells.Item(39, 6).Value = 7
$ws.Cells.Item(40, 6).Value = 30
$ws.Cells.Item(42, 6).Value = 24
$ws.Cells.Item(43, 6).Value = 24
$ws.Cells.Item(45, 6).Value = 109
$ws.Cells.Item(47, 6).Value = 3

# --- Sheet 3: 本地生活 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(4, 6).Value = 860
$ws.Cells.Item(6, 6).Value = 2667
$ws.Cells.Item(7, 6).Value = 4442
$ws.Cells.Item(8, 6).Value = 111
$ws.Cells.Item(10, 6).Value = 496
$ws.Cells.Item(11, 6).Value = 569
$ws.Cells.Item(12, 6).Value = 382
$ws.Cells.Item(13, 6).Value = 111
$ws.Cells.Item(14, 6).Value = 535
$ws.Cells.Item(15, 6).Value = 165
$ws.Cells.Item(16, 6).Value = 212

# --- Sheet 4: 全部类型 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 1826
$ws.Cells.Item(3, 6).Value = 860
$ws.Cells.Item(5, 6).Value = 4442
$ws.Cells.Item(6, 6).Value = 111
$ws.Cells.Item(7, 6).Value = 569
$ws.Cells.Item(8, 6).Value = 34
$ws.Cells.Item(9, 6).Value = 111
$ws.Cells.Item(10, 6).Value = 111
$ws.Cells.Item(11, 6).Value = 535
$ws.Cells.Item(12, 6).Value = 165
$ws.Cells.Item(14, 6).Value = 1076
$ws.Cells.Item(15, 6).Value = 40
$ws.Cells.Item(16, 6).Value = 167
$ws.Cells.Item(18, 6).Value = 1371
$ws.Cells.Item(19, 6).Value = 566
$ws.Cells.Item(20, 6).Value = 444
$ws.Cells.Item(21, 6).Value = 204
$ws.Cells.Item(22, 6).Value = 2460
$ws.Cells.Item(23, 6).Value = 1164
$ws.Cells.Item(24, 6).Value = 1358
$ws.Cells.Item(25, 6).Value = 1193
$ws.Cells.Item(26, 6).Value = 1395
$ws.Cells.Item(27, 6).Value = 109
$ws.Cells.Item(28, 6).Value = 199
$ws.Cells.Item(29, 6).Value = 59
$ws.Cells.Item(30, 6).Value = 1539
$ws.Cells.Item(31, 6).Value = 763
$ws.Cells.Item(32, 6).Value = 307
$ws.Cells.Item(33, 6).Value = 414
$ws.Cells.Item(34, 6).Value = 1167
$ws.Cells.Item(36, 6).Value = 541
$ws.Cells.Item(37, 6).Value = 995
$ws.Cells.Item(38, 6).Value = 267
$ws.Cells.Item(40, 6).Value = 875
$ws.Cells.Item(41, 6).Value = 815
$ws.Cells.Item(43, 6).Value = 1539
$ws.Cells.Item(44, 6).Value = 83
$ws.Cells.Item(45, 6).Value = 161
$ws.Cells.Item(47, 6).Value = 774
$ws.Cells.Item(48, 6).Value = 24
$ws.Cells.Item(49, 6).Value = 758
$ws.Cells.Item(50, 6).Value = 100
